$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SF64057"
$ws.Range("A3").Value = "SF289805"
$ws.Range("A4").Value = "SF128322"
$ws.Range("A5").Value = "SF26908"
$ws.Range("A6").Value = "SF73266"
